# In the schedule table, the leftmost ("Заняття" / class date) column of
# two consecutive rows is still blank - the rows identified by their
# lesson-code cell ("Л04" and "ПР02"). Word originally left these
# paragraphs empty (only pPr/rPr, no run). Fill them in with the next
# two class dates in the sequence, "27.09" and "28.09", using the same
# Times New Roman / 14pt formatting already carried by the paragraph.

$d = $word.ActiveDocument

function Get-DateCell($markerText) {
    $rng = $d.Content
    $rng.Find.Execute($markerText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $row = $rng.Cells.Item(1).Row
    return $row.Cells.Item(1)
}

function Set-DateCell($markerText, $dateText) {
    $cell = Get-DateCell $markerText
    $cell.Range.Text = $dateText

    # Re-locate the same cell before every subsequent formatting call -
    # each package mutation can make the previous handle stale.
    $cell = Get-DateCell $markerText
    $cell.Range.Font.Name = "Times New Roman"

    $cell = Get-DateCell $markerText
    $cell.Range.Font.Size = 14

    $cell = Get-DateCell $markerText
    $cell.Range.Font.SizeBi = 14
}

Set-DateCell "Л04" "27.09"
Set-DateCell "ПР02" "28.09"
